# Regenerate all Word files with proper table formatting
#
# 1. Remove the old front-matter navigation paragraphs ("Home",
#    "<- Back to Home", "Download Word Document") that used to sit
#    before the document heading.
# 2. Give every table an explicit 100% preferred width
#    (<w:tblW w:type="pct" w:w="5000"/>) instead of "auto".

$d = $word.ActiveDocument

# --- 1. Drop the three leading navigation paragraphs -----------------
# They are always the first three paragraphs of the body: "Home",
# "<- Back to Home" and "Download Word Document", each its own
# hyperlinked paragraph right before the "Team Meeting Agenda..." H1.
# Guard on the visible text so we only ever remove that exact block.
$p1 = $d.Paragraphs(1).Range.Text.Trim()
$p2 = $d.Paragraphs(2).Range.Text.Trim()
$p3 = $d.Paragraphs(3).Range.Text.Trim()

if ($p1 -eq "Home" -and $p2.EndsWith("Back to Home") -and $p3 -eq "Download Word Document") {
    $navStart = $d.Paragraphs(1).Range.Start
    $navEnd = $d.Paragraphs(3).Range.End
    $navRange = $d.Range($navStart, $navEnd)
    $navRange.Delete()
}

# --- 2. Force every table to 100% preferred width ---------------------
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    $tbl.PreferredWidthType = 2   # wdPreferredWidthPercent
    $tbl.PreferredWidth = 250     # 250 * 20 = 5000 (100% in fiftieths-of-a-percent)
}
